$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (G1) onto the
# new "Save" header cell (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" data column values for the two data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
